$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry below is (cell address, new value) reflecting the updated
# crypto price/volume/hour snapshot for this run (symbol list refresh).
# Text number format ("@") is applied before writing so that numeric-
# looking strings (prices with trailing zeros, percentages, small
# decimals like 0.00002101) are stored exactly as text, matching the
# original inline-string cell contents instead of being reinterpreted
# as numbers/percentages by Excel.
$updates = @(
    @{ Cell = "D2"; Value = "321.55" }
    @{ Cell = "E2"; Value = "2.88%" }
    @{ Cell = "G2"; Value = "3" }
    @{ Cell = "D3"; Value = "39.81" }
    @{ Cell = "E3"; Value = "5.43%" }
    @{ Cell = "G3"; Value = "3" }
    @{ Cell = "D4"; Value = "5.239" }
    @{ Cell = "E4"; Value = "2.09%" }
    @{ Cell = "G4"; Value = "3" }
    @{ Cell = "D5"; Value = "0.08103" }
    @{ Cell = "E5"; Value = "2.49%" }
    @{ Cell = "G5"; Value = "3" }
    @{ Cell = "D6"; Value = "4.519" }
    @{ Cell = "E6"; Value = "2.30%" }
    @{ Cell = "G6"; Value = "3" }
    @{ Cell = "D7"; Value = "8.602" }
    @{ Cell = "E7"; Value = "3.83%" }
    @{ Cell = "G7"; Value = "3" }
    @{ Cell = "D8"; Value = "1.918" }
    @{ Cell = "E8"; Value = "0.67%" }
    @{ Cell = "G8"; Value = "3" }
    @{ Cell = "D9"; Value = "2.942" }
    @{ Cell = "E9"; Value = "-1.55%" }
    @{ Cell = "G9"; Value = "3" }
    @{ Cell = "D10"; Value = "0.9366" }
    @{ Cell = "E10"; Value = "1.48%" }
    @{ Cell = "G10"; Value = "3" }
    @{ Cell = "D11"; Value = "0.1290" }
    @{ Cell = "E11"; Value = "13.21%" }
    @{ Cell = "G11"; Value = "3" }
    @{ Cell = "D12"; Value = "0.1958" }
    @{ Cell = "E12"; Value = "3.72%" }
    @{ Cell = "G12"; Value = "3" }
    @{ Cell = "D13"; Value = "0.09182" }
    @{ Cell = "E13"; Value = "0.35%" }
    @{ Cell = "G13"; Value = "3" }
    @{ Cell = "D14"; Value = "0.03431" }
    @{ Cell = "E14"; Value = "3.15%" }
    @{ Cell = "G14"; Value = "3" }
    @{ Cell = "D15"; Value = "0.09542" }
    @{ Cell = "E15"; Value = "-0.68%" }
    @{ Cell = "G15"; Value = "3" }
    @{ Cell = "D16"; Value = "0.001405" }
    @{ Cell = "E16"; Value = "1.95%" }
    @{ Cell = "G16"; Value = "3" }
    @{ Cell = "B17"; Value = "TigerCash" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" }
    @{ Cell = "D17"; Value = "0.006464" }
    @{ Cell = "E17"; Value = "4.70%" }
    @{ Cell = "G17"; Value = "3" }
    @{ Cell = "B18"; Value = "LEO" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" }
    @{ Cell = "D18"; Value = "3.357" }
    @{ Cell = "E18"; Value = "-5.66%" }
    @{ Cell = "G18"; Value = "3" }
    @{ Cell = "B19"; Value = "BitpandaEcosystemToken" }
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best" }
    @{ Cell = "D19"; Value = "0.3535" }
    @{ Cell = "E19"; Value = "2.55%" }
    @{ Cell = "G19"; Value = "3" }
    @{ Cell = "B20"; Value = "MCDex" }
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" }
    @{ Cell = "D20"; Value = "6.748" }
    @{ Cell = "E20"; Value = "27.50%" }
    @{ Cell = "G20"; Value = "3" }
    @{ Cell = "B21"; Value = "ProBitToken" }
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob" }
    @{ Cell = "D21"; Value = "0.1328" }
    @{ Cell = "E21"; Value = "3.12%" }
    @{ Cell = "G21"; Value = "3" }
    @{ Cell = "B22"; Value = "ZBToken" }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb" }
    @{ Cell = "D22"; Value = "0.2312" }
    @{ Cell = "E22"; Value = "-10.71%" }
    @{ Cell = "G22"; Value = "3" }
    @{ Cell = "B23"; Value = "CoinExToken" }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" }
    @{ Cell = "D23"; Value = "0.04438" }
    @{ Cell = "E23"; Value = "1.71%" }
    @{ Cell = "G23"; Value = "3" }
    @{ Cell = "D24"; Value = "0.001223" }
    @{ Cell = "E24"; Value = "-0.79%" }
    @{ Cell = "G24"; Value = "3" }
    @{ Cell = "D25"; Value = "0.004353" }
    @{ Cell = "G25"; Value = "3" }
    @{ Cell = "E26"; Value = "-5.13%" }
    @{ Cell = "G26"; Value = "3" }
    @{ Cell = "D27"; Value = "0.0003991" }
    @{ Cell = "E27"; Value = "0.00%" }
    @{ Cell = "G27"; Value = "3" }
    @{ Cell = "G28"; Value = "3" }
    @{ Cell = "G29"; Value = "3" }
    @{ Cell = "G30"; Value = "3" }
    @{ Cell = "G31"; Value = "3" }
    @{ Cell = "G32"; Value = "3" }
    @{ Cell = "G33"; Value = "3" }
    @{ Cell = "G34"; Value = "3" }
    @{ Cell = "G35"; Value = "3" }
    @{ Cell = "G36"; Value = "3" }
    @{ Cell = "G37"; Value = "3" }
    @{ Cell = "G38"; Value = "3" }
    @{ Cell = "D39"; Value = "0.02438" }
    @{ Cell = "E39"; Value = "7.56%" }
    @{ Cell = "G39"; Value = "3" }
    @{ Cell = "D40"; Value = "0.05206" }
    @{ Cell = "E40"; Value = "2.21%" }
    @{ Cell = "G40"; Value = "3" }
    @{ Cell = "D41"; Value = "0.007694" }
    @{ Cell = "E41"; Value = "3.25%" }
    @{ Cell = "G41"; Value = "3" }
    @{ Cell = "D42"; Value = "0.1430" }
    @{ Cell = "E42"; Value = "5.57%" }
    @{ Cell = "G42"; Value = "3" }
    @{ Cell = "D43"; Value = "0.008796" }
    @{ Cell = "E43"; Value = "-2.33%" }
    @{ Cell = "G43"; Value = "3" }
    @{ Cell = "E44"; Value = "5.52%" }
    @{ Cell = "G44"; Value = "3" }
    @{ Cell = "D45"; Value = "0.008163" }
    @{ Cell = "E45"; Value = "-5.63%" }
    @{ Cell = "G45"; Value = "3" }
    @{ Cell = "E46"; Value = "-0.56%" }
    @{ Cell = "G46"; Value = "3" }
    @{ Cell = "E47"; Value = "0.00%" }
    @{ Cell = "G47"; Value = "3" }
    @{ Cell = "D48"; Value = "0.002852" }
    @{ Cell = "E48"; Value = "-12.22%" }
    @{ Cell = "G48"; Value = "3" }
    @{ Cell = "D49"; Value = "0.002481" }
    @{ Cell = "E49"; Value = "148.00%" }
    @{ Cell = "G49"; Value = "3" }
    @{ Cell = "D50"; Value = "0.00002101" }
    @{ Cell = "E50"; Value = "0.00%" }
    @{ Cell = "G50"; Value = "3" }
    @{ Cell = "D51"; Value = "0.0002001" }
    @{ Cell = "E51"; Value = "0.00%" }
    @{ Cell = "G51"; Value = "3" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
